$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.7955796666666667
$ws.Range("H2").Value = 2.386739
$ws.Range("I2").Value = 0.1186174580157865
$ws.Range("J2").Value = 0.1186174580157865
$ws.Range("M2").Value = 9.771369666666667
$ws.Range("N2").Value = 29.314109
$ws.Range("O2").Value = 0.6454156383975566
$ws.Range("P2").Value = 0.6454156383975566
$ws.Range("Q2").Value = 7.773903022283445
$ws.Range("R2").Value = 69.96512720055101
$ws.Range("S2").Value = 0.07655756239035423
$ws.Range("T2").Value = 0.07655756239035423
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.7955796666666667
$ws.Range("H3").Value = 2.386739
$ws.Range("I3").Value = 0.1186174580157865
$ws.Range("J3").Value = 0.1186174580157865
$ws.Range("O3").Value = 0.1821792144395723
$ws.Range("P3").Value = 0.1821792144395723
$ws.Range("Q3").Value = 2.194312411216556
$ws.Range("R3").Value = 19.748811700949
$ws.Range("S3").Value = 0.02160963532013494
$ws.Range("T3").Value = 0.02160963532013493
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.7955796666666667
$ws.Range("H4").Value = 2.386739
$ws.Range("I4").Value = 0.1186174580157865
$ws.Range("J4").Value = 0.1186174580157865
$ws.Range("M4").Value = 1.889356
$ws.Range("N4").Value = 5.668068
$ws.Range("O4").Value = 0.1247951874198449
$ws.Range("P4").Value = 0.1247951874198449
$ws.Range("Q4").Value = 1.503133216694667
$ws.Range("R4").Value = 13.528198950252
$ws.Range("S4").Value = 0.01480288790434566
$ws.Range("T4").Value = 0.01480288790434566
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.7955796666666667
$ws.Range("H5").Value = 2.386739
$ws.Range("I5").Value = 0.1186174580157865
$ws.Range("J5").Value = 0.1186174580157865
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.7207983333333333
$ws.Range("N5").Value = 2.162395
$ws.Range("O5").Value = 0.04760995974302628
$ws.Range("P5").Value = 0.04760995974302627
$ws.Range("Q5").Value = 0.5734524977672222
$ws.Range("R5").Value = 5.161072479905
$ws.Range("S5").Value = 0.005647372400951707
$ws.Range("T5").Value = 0.005647372400951706
$ws.Range("I6").Value = 0.6312226244877757
$ws.Range("J6").Value = 0.6312226244877758
$ws.Range("M6").Value = 9.771369666666667
$ws.Range("N6").Value = 29.314109
$ws.Range("O6").Value = 0.6454156383975566
$ws.Range("P6").Value = 0.6454156383975566
$ws.Range("Q6").Value = 41.36881324489467
$ws.Range("R6").Value = 372.3193192040521
$ws.Range("S6").Value = 0.4074009531547588
$ws.Range("T6").Value = 0.4074009531547589
$ws.Range("I7").Value = 0.6312226244877757
$ws.Range("J7").Value = 0.6312226244877758
$ws.Range("O7").Value = 0.1821792144395723
$ws.Range("P7").Value = 0.1821792144395723
$ws.Range("S7").Value = 0.1149956418656681
$ws.Range("T7").Value = 0.1149956418656681
$ws.Range("I8").Value = 0.6312226244877757
$ws.Range("J8").Value = 0.6312226244877758
$ws.Range("M8").Value = 1.889356
$ws.Range("N8").Value = 5.668068
$ws.Range("O8").Value = 0.1247951874198449
$ws.Range("P8").Value = 0.1247951874198449
$ws.Range("Q8").Value = 7.998921152656
$ws.Range("R8").Value = 71.990290373904
$ws.Range("S8").Value = 0.07877354572659834
$ws.Range("T8").Value = 0.07877354572659832
$ws.Range("I9").Value = 0.6312226244877757
$ws.Range("J9").Value = 0.6312226244877758
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.7207983333333333
$ws.Range("N9").Value = 2.162395
$ws.Range("O9").Value = 0.04760995974302628
$ws.Range("P9").Value = 0.04760995974302627
$ws.Range("Q9").Value = 3.051626604673333
$ws.Range("R9").Value = 27.46463944206
$ws.Range("S9").Value = 0.03005248374075039
$ws.Range("T9").Value = 0.03005248374075039
$ws.Range("G10").Value = 1.666370333333333
$ws.Range("H10").Value = 4.999111
$ws.Range("I10").Value = 0.248448548064433
$ws.Range("J10").Value = 0.248448548064433
$ws.Range("M10").Value = 9.771369666666667
$ws.Range("N10").Value = 29.314109
$ws.Range("O10").Value = 0.6454156383975566
$ws.Range("P10").Value = 0.6454156383975566
$ws.Range("Q10").Value = 16.28272052856655
$ws.Range("R10").Value = 146.544484757099
$ws.Range("S10").Value = 0.160352578257952
$ws.Range("T10").Value = 0.160352578257952
$ws.Range("G11").Value = 1.666370333333333
$ws.Range("H11").Value = 4.999111
$ws.Range("I11").Value = 0.248448548064433
$ws.Range("J11").Value = 0.248448548064433
$ws.Range("O11").Value = 0.1821792144395723
$ws.Range("P11").Value = 0.1821792144395723
$ws.Range("Q11").Value = 4.596066562933444
$ws.Range("R11").Value = 41.364599066401
$ws.Range("S11").Value = 0.04526216131503072
$ws.Range("T11").Value = 0.04526216131503071
$ws.Range("G12").Value = 1.666370333333333
$ws.Range("H12").Value = 4.999111
$ws.Range("I12").Value = 0.248448548064433
$ws.Range("J12").Value = 0.248448548064433
$ws.Range("M12").Value = 1.889356
$ws.Range("N12").Value = 5.668068
$ws.Range("O12").Value = 0.1247951874198449
$ws.Range("P12").Value = 0.1247951874198449
$ws.Range("Q12").Value = 3.148366787505333
$ws.Range("R12").Value = 28.335301087548
$ws.Range("S12").Value = 0.03100518311988925
$ws.Range("T12").Value = 0.03100518311988925
$ws.Range("G13").Value = 1.666370333333333
$ws.Range("H13").Value = 4.999111
$ws.Range("I13").Value = 0.248448548064433
$ws.Range("J13").Value = 0.248448548064433
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.7207983333333333
$ws.Range("N13").Value = 2.162395
$ws.Range("O13").Value = 0.04760995974302628
$ws.Range("P13").Value = 0.04760995974302627
$ws.Range("Q13").Value = 1.201116958982778
$ws.Range("R13").Value = 10.810052630845
$ws.Range("S13").Value = 0.01182862537156098
$ws.Range("T13").Value = 0.01182862537156098
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.01147833333333333
$ws.Range("H14").Value = 0.034435
$ws.Range("I14").Value = 0.00171136943200476
$ws.Range("J14").Value = 0.00171136943200476
$ws.Range("M14").Value = 9.771369666666667
$ws.Range("N14").Value = 29.314109
$ws.Range("O14").Value = 0.6454156383975566
$ws.Range("P14").Value = 0.6454156383975566
$ws.Range("Q14").Value = 0.1121590381572222
$ws.Range("R14").Value = 1.009431343415
$ws.Range("S14").Value = 0.001104544594491416
$ws.Range("T14").Value = 0.001104544594491416
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.01147833333333333
$ws.Range("H15").Value = 0.034435
$ws.Range("I15").Value = 0.00171136943200476
$ws.Range("J15").Value = 0.00171136943200476
$ws.Range("O15").Value = 0.1821792144395723
$ws.Range("P15").Value = 0.1821792144395723
$ws.Range("Q15").Value = 0.03165873934277778
$ws.Range("R15").Value = 0.284928654085
$ws.Range("S15").Value = 0.0003117759387385243
$ws.Range("T15").Value = 0.0003117759387385242
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.01147833333333333
$ws.Range("H16").Value = 0.034435
$ws.Range("I16").Value = 0.00171136943200476
$ws.Range("J16").Value = 0.00171136943200476
$ws.Range("M16").Value = 1.889356
$ws.Range("N16").Value = 5.668068
$ws.Range("O16").Value = 0.1247951874198449
$ws.Range("P16").Value = 0.1247951874198449
$ws.Range("Q16").Value = 0.02168665795333333
$ws.Range("R16").Value = 0.19517992158
$ws.Range("S16").Value = 0.0002135706690116275
$ws.Range("T16").Value = 0.0002135706690116275
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.01147833333333333
$ws.Range("H17").Value = 0.034435
$ws.Range("I17").Value = 0.00171136943200476
$ws.Range("J17").Value = 0.00171136943200476
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.7207983333333333
$ws.Range("N17").Value = 2.162395
$ws.Range("O17").Value = 0.04760995974302628
$ws.Range("P17").Value = 0.04760995974302627
$ws.Range("Q17").Value = 0.008273563536111112
$ws.Range("R17").Value = 0.074462071825
$ws.Range("S17").Value = 0.00008147822976319238
$ws.Range("T17").Value = 0.00008147822976319237
